$d = $word.ActiveDocument

# --- Change 1: SSU Lobby "Start the game" button description ---------------
# Original: takođe ima dugme ”Start the game” kojim se prelazi u prikaz Igra i
#           vrši shodna funkcionalnost za svakog igrača trenutno u lobby-u,
#           pored dugmeta ”Exit”.
# New:      takođe ima dugme ”Start the game” , pored dugmeta ”Exit”, kojim on
#           pokreće igru, čime se prelazi se u prikaz Igra i vrši shodna
#           funkcionalnost za svakog igrača trenutno u lobby-u.

$r = $d.Content
$r.Find.Execute(" kojim se prelazi u prikaz Igra i vrši shodna funkcionalnost", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = " , pored dugmeta ”Exit”, kojim on pokreće igru, čime se prelazi se u prikaz Igra i vrši shodna funkcionalnost"

$r2 = $d.Content
$r2.Find.Execute(", pored dugmeta ”Exit”.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Text = "."

# --- Change 2: host leaves the lobby --------------------------------------
# Original: Ukoliko domaćin igre odluči da izađe iz lobby-a u koraku 1 ...
# New:      Ukoliko domaćin lobby-a odluči da izađe iz lobby-a u koraku 1 ...

$r3 = $d.Content
$r3.Find.Execute("domaćin igre odluči da izađe iz", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Text = "domaćin lobby-a odluči da izađe iz"
